# Update "想去人数" (want-to-go count) figures in column F for a handful of
# events on both the "展览" sheet and the aggregated "全部类型" sheet, to
# reflect newly generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F16").Value = 6568
$wsExhibit.Range("F20").Value = 169
$wsExhibit.Range("F23").Value = 15701
$wsExhibit.Range("F26").Value = 302
$wsExhibit.Range("F32").Value = 268
$wsExhibit.Range("F35").Value = 312

# Sheet "全部类型" (All Types) - same events, different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F19").Value = 6568
$wsAll.Range("F23").Value = 169
$wsAll.Range("F27").Value = 15701
$wsAll.Range("F30").Value = 302
$wsAll.Range("F37").Value = 268
$wsAll.Range("F40").Value = 312
